$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Fix mis-encoded (mojibake) Spanish accented text in the "Zona",
# "Calle" and "Municipio" columns (B, C, D) of the sensor table.
# ------------------------------------------------------------------

# Row 6 - ANL13 / Juarez
$ws.Range("B6").Value = "Juarez"
$ws.Range("D6").Value = "Ju" + [char]0x00E1 + "rez"

# Row 7 - ANL15 / Pesqueria
$ws.Range("C7").Value = "Boulevard Rogelio A. P" + [char]0x00E9 + "rez Arrambide"
$ws.Range("D7").Value = "Pesquer" + [char]0x00ED + "a"

# Row 8 - ANL16 / San Juan
$ws.Range("C8").Value = "Av. Misi" + [char]0x00F3 + "n Arcos 69, Misi" + [char]0x00F3 + "n San Juan"
$ws.Range("D8").Value = "Garc" + [char]0x00ED + "a"

# Row 9 - ANL2 / San Nicolas
$ws.Range("B9").Value = "San Nicol" + [char]0x00E1 + "s"
$ws.Range("D9").Value = "San Nicol" + [char]0x00E1 + "s de los Garza"

# Row 11 - ANL4 / San Pedro
$ws.Range("C11").Value = "GRAL. GARZA AYALA ESQUINA CON DIEGO SALD" + [char]0x00CD + "VAR"

# Row 12 - ANL5 / San Nicolas (UANL)
$ws.Range("B12").Value = "San Nicol" + [char]0x00E1 + "s (UANL)"
$ws.Range("D12").Value = "San Nicol" + [char]0x00E1 + "s de los Garza"

# Row 13 - ANL6 / Garcia
$ws.Range("B13").Value = "Garc" + [char]0x00ED + "a"
$ws.Range("D13").Value = "Garc" + [char]0x00ED + "a"

# Row 14 - ANL7 / San Bernabe
$ws.Range("B14").Value = "San Bernab" + [char]0x00E9

# Row 15 - ANL8 / Cadereyta
$ws.Range("D15").Value = "Cadereyta Jim" + [char]0x00E9 + "nez"

# ------------------------------------------------------------------
# Column widths for B, C, D
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 13.75
$ws.Columns.Item(3).ColumnWidth = 47.42
$ws.Columns.Item(4).ColumnWidth = 28.09

# ------------------------------------------------------------------
# View state: zoom + selection
# ------------------------------------------------------------------
[void]$ws.Range("C19").Select()
$excel.ActiveWindow.Zoom = 150
